# Generate Report for Handback
#
# For each locale sheet (zh-cn, de-de):
#   - Status moves from "Ready for handoff" -> "Handed back: in sync with en-US"
#   - Two new columns are populated for rows 2 & 3: "Latest Target File" (E) and
#     "Latest Handback File" (F), each holding a hyperlinked file name (mirrors the
#     Source File/Latest Handoff File values, which is what the source report does).
#   - "Latest Handback DateTime" (G) moves from the placeholder 0001-01-01 00:00:00 to
#     a real handback timestamp.

$wb = $excel.ActiveWorkbook

# Hyperlink target URLs reused from the existing links on each sheet.
$aMdUrl     = "https://github.com/OpenLocalizationTest/oltest/blob/766bca4b26c2296daff10bd35fb505dfc891b906/e2e/a.md.md"
$zhCnXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b0f95d26428b7281ae218944527a9d1c2608f93c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf"
$deDeXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f34b63d120c81de71c85aace8315ee6445bb1157/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf"

# Blue underlined "hyperlink" look used by the other linked cells (A2, C2, ...).
$hyperlinkColor = 15570276   # BGR int for RGB FF6495ED

function Apply-HandbackRow($ws, $row, $targetFile, $targetFileUrl, $handbackFile, $handbackFileUrl, $handbackDateTime) {
    # Status -> handed back, in sync with en-US
    $ws.Cells.Item($row, 2).Value = "Handed back: in sync with en-US"

    # E: Latest Target File (new hyperlinked cell)
    $eCell = $ws.Cells.Item($row, 5)
    $eCell.Value = $targetFile
    $ws.Hyperlinks.Add($eCell, $targetFileUrl, [Type]::Missing, [Type]::Missing, $targetFile) | Out-Null
    $eCell.Font.Underline = 2
    $eCell.Font.Color = $hyperlinkColor

    # F: Latest Handback File (new hyperlinked cell)
    $fCell = $ws.Cells.Item($row, 6)
    $fCell.Value = $handbackFile
    $ws.Hyperlinks.Add($fCell, $handbackFileUrl, [Type]::Missing, [Type]::Missing, $handbackFile) | Out-Null
    $fCell.Font.Underline = 2
    $fCell.Font.Color = $hyperlinkColor

    # G: Latest Handback DateTime -> real timestamp instead of the 0001-01-01 placeholder
    $ws.Cells.Item($row, 7).Value = $handbackDateTime
}

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
Apply-HandbackRow $wsZh 2 "a.md.md" $aMdUrl "a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf" $zhCnXlfUrl "2016-01-25 06:53:41"
Apply-HandbackRow $wsZh 3 "a.md.md" $aMdUrl "a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf" $zhCnXlfUrl "2016-01-25 06:53:41"

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
Apply-HandbackRow $wsDe 2 "a.md.md" $aMdUrl "a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf" $deDeXlfUrl "2016-01-25 06:53:59"
Apply-HandbackRow $wsDe 3 "a.md.md" $aMdUrl "a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf" $deDeXlfUrl "2016-01-25 06:53:59"

Write-Host "Handback report generated."
